# HW_2/Summary.xlsx edit:
# 1. Add a new result row for "Proximity Search (Unmodified queries)" under the
#    STEMMED / "With HEAD - TEXT Indexed EC2" section (new row 13), pushing the
#    rest of the sheet (NON STEMMED section and below) down by one row.
# 2. Leave a blank separator row (row 14) before the NON STEMMED section, same
#    as the blank separator rows used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 - shifts everything from the old row 13
# (the blank separator before "NON STEMMED") downward by one.
$ws.Rows.Item(13).Insert()

# Populate the new row with the proximity-search-with-unmodified-queries result.
$ws.Range("A13").Value2 = "Proximity Search (Unmodified queries)"
$ws.Range("B13").Value2 = 0.2194
$ws.Range("C13").Value2 = 0.34
$ws.Range("D13").Value2 = 0.2827

# Match the author's final selection/active cell.
$ws.Range("A13").Select()
